# Update the "OpsTracker" sheet of the Daily Status Tracker workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpsTracker")

# Re-assign ownership of a few open items from "Victor" to "Rahul".
$ws.Range("C32").Value = "Rahul"
$ws.Range("C33").Value = "Rahul"
$ws.Range("C34").Value = "Rahul"

# Add follow-up comments for the two "File" related items.
$ws.Range("E32").Value = "New File to be purchased"
$ws.Range("E33").Value = "Already file is there, we need to remove other documents from this file"

# Update the saved selection to reflect where the user left off editing.
$ws.Activate()
$ws.Range("E37").Select()
